$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.068.14"
$ws.Range("E2").Value = "  -1.27%  "

# Row 3
$ws.Range("D3").Value = "1.830.36"
$ws.Range("E3").Value = "  -1.33%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.41%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.42%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2946"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.73%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07335"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.57%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "

# Row 12
$ws.Range("D12").Value = "1.830.94"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.007"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6722"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.60%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.74%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.111"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.82%  "

# Row 17
$ws.Range("D17").Value = "29.045.04"
$ws.Range("E17").Value = "  -1.29%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008211"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.255"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.651"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.67%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.62%  "

# Row 28
$ws.Range("E28").Value = "  -2.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.222"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30
$ws.Range("E30").Value = "  -1.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.196"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.26%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05320"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7476"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.680"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("D37").Value = "1.296.18"
$ws.Range("E37").Value = "  -2.57%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01805"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.702"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9225"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.40%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.73%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9985"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.50%  "

# Row 44
$ws.Range("D44").Value = "1.976.53"
$ws.Range("E44").Value = "  -1.26%  "

# Row 45
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5175"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46
$ws.Range("B46").Value = "XinFinNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07690"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +16.02%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.88%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.748"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.66%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.235"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05906"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.27%  "
